{"js": "const pairs = [\n  [\"2025-08-23 Saturday\", \"2025-08-24 Sunday\"],\n  [\"508\u00f78=63, 4\", \"342\u00f72=171, 0\"],\n  [\"469\u00f75=93, 4\", \"859\u00f73=286, 1\"],\n  [\"446\u00f78=55, 6\", \"762\u00f73=254, 0\"],\n  [\"259\u00f77=37, 0\", \"682\u00f78=85, 2\"],\n  [\"505\u00f74=126, 1\", \"705\u00f75=141, 0\"],\n  [\"652\u00f77=93, 1\", \"631\u00f77=90, 1\"],\n  [\"877\u00f73=292, 1\", \"585\u00f78=73, 1\"],\n  [\"926\u00f76=154, 2\", \"428\u00f79=47, 5\"],\n  [\"785\u00f77=112, 1\", \"728\u00f74=182, 0\"],\n  [\"644\u00f76=107, 2\", \"869\u00f72=434, 1\"],\n  [\"476\u00f73=158, 2\", \"182\u00f72=91, 0\"],\n  [\"766\u00f75=153, 1\", \"156\u00f77=22, 2\"],\n  [\"288\u00f74=72, 0\", \"340\u00f79=37, 7\"],\n  [\"843\u00f74=210, 3\", \"185\u00f76=30, 5\"],\n  [\"498\u00f72=249, 0\", \"539\u00f75=107, 4\"],\n  [\"220\u00f74=55, 0\", \"574\u00f79=63, 7\"],\n  [\"323\u00f72=161, 1\", \"145\u00f76=24, 1\"],\n  [\"821\u00f78=102, 5\", \"622\u00f72=311, 0\"],\n  [\"182\u00f74=45, 2\", \"472\u00f74=118, 0\"],\n  [\"698\u00f74=174, 2\", \"458\u00f72=229, 0\"],\n  [\"505\u00f77=72, 1\", \"771\u00f76=128, 3\"],\n  [\"948\u00f74=237, 0\", \"711\u00f76=118, 3\"],\n  [\"925\u00f79=102, 7\", \"941\u00f74=235, 1\"],\n  [\"632\u00f74=158, 0\", \"973\u00f72=486, 1\"],\n  [\"918\u00f73=306, 0\", \"607\u00f73=202, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of pairs) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + find);\n  }\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$wdReplaceAll = 2\n$wdFindContinue = 1\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2025-08-23 Saturday', '2025-08-24 Sunday'),\n    @('508\u00f78=63, 4', '342\u00f72=171, 0'),\n    @('469\u00f75=93, 4', '859\u00f73=286, 1'),\n    @('446\u00f78=55, 6', '762\u00f73=254, 0'),\n    @('259\u00f77=37, 0', '682\u00f78=85, 2'),\n    @('505\u00f74=126, 1', '705\u00f75=141, 0'),\n    @('652\u00f77=93, 1', '631\u00f77=90, 1'),\n    @('877\u00f73=292, 1', '585\u00f78=73, 1'),\n    @('926\u00f76=154, 2', '428\u00f79=47, 5'),\n    @('785\u00f77=112, 1', '728\u00f74=182, 0'),\n    @('644\u00f76=107, 2', '869\u00f72=434, 1'),\n    @('476\u00f73=158, 2', '182\u00f72=91, 0'),\n    @('766\u00f75=153, 1', '156\u00f77=22, 2'),\n    @('288\u00f74=72, 0', '340\u00f79=37, 7'),\n    @('843\u00f74=210, 3', '185\u00f76=30, 5'),\n    @('498\u00f72=249, 0', '539\u00f75=107, 4'),\n    @('220\u00f74=55, 0', '574\u00f79=63, 7'),\n    @('323\u00f72=161, 1', '145\u00f76=24, 1'),\n    @('821\u00f78=102, 5', '622\u00f72=311, 0'),\n    @('182\u00f74=45, 2', '472\u00f74=118, 0'),\n    @('698\u00f74=174, 2', '458\u00f72=229, 0'),\n    @('505\u00f77=72, 1', '771\u00f76=128, 3'),\n    @('948\u00f74=237, 0', '711\u00f76=118, 3'),\n    @('925\u00f79=102, 7', '941\u00f74=235, 1'),\n    @('632\u00f74=158, 0', '973\u00f72=486, 1'),\n    @('918\u00f73=306, 0', '607\u00f73=202, 1'),\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Replace failed for: $findText\"\n    }\n}\n"}
